$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# Insert two new personnel rows above the existing row 6 (Katherine Qi),
# shifting rows 6-7 down to 8-9.
$ws.Range("A6:A7").EntireRow.Insert()

# New row 6: E. Taylor Crockford, technician
$ws.Range("A6").Value = "E. Taylor"
$ws.Range("C6").Value = "Crockford"
$ws.Range("D6").Value = "Northeast U.S. Shelf LTER"
$ws.Range("E6").Value = "ecrockford@whoi.edu"
$ws.Range("G6").Value = "technician"
$ws.Range("H6").Value = "Northeast U.S. Shelf LTER"
$ws.Range("I6").Value = "NSF"
$ws.Range("J6").Value = "OCE-1655686"

# New row 7: Joe Futrelle, softwareDeveloper
$ws.Range("A7").Value = "Joe"
$ws.Range("C7").Value = "Futrelle"
$ws.Range("D7").Value = "Northeast U.S. Shelf LTER"
$ws.Range("E7").Value = "jfutrelle@whoi.edu"
$ws.Range("G7").Value = "softwareDeveloper"
$ws.Range("H7").Value = "Northeast U.S. Shelf LTER"
$ws.Range("I7").Value = "NSF"
$ws.Range("J7").Value = "OCE-1655686"

# Make the Personnel sheet the active/selected tab (it previously was
# Keywords); this also clears tabSelected on the Keywords sheet and
# resets the workbook's active tab.
$ws.Activate()
